$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete column H ("Rating2" column) - this removes the Rating2 header/formula
#    cells and shifts I:M left to H:L (widths, text, merges, dimension all follow).
$ws.Columns("H").Delete()

# 2. Row 3 used to be a second copy of the header labels; it should instead hold the
#    "&=result.X" formula placeholders (one column to the left of where they used to
#    live, now that the Rating2 column is gone). Rewrite it explicitly.
$ws.Range("A3").Value2 = "&=result.AuditDate"
$ws.Range("B3").Value2 = "&=result.AuditType"
$ws.Range("C3").Value2 = "&=result.AuditType"
$ws.Range("D3").Value2 = "&=result.LineID"
$ws.Range("E3").Value2 = "&=result.RatingNa"
$ws.Range("F3").Value2 = "&=result.Rating0"
$ws.Range("G3").Value2 = "&=result.Rating1"
$ws.Range("H3").Value2 = "&=result.Total"
$ws.Range("I3").Value2 = "&=result.NeedToDoQty"
$ws.Range("J3").Value2 = "&=result.Achieving"

# 3. A3 (Audit Date placeholder) now carries a date format + wrapped text.
$ws.Range("A3").NumberFormat = "yyyy/mm/dd"
$ws.Range("A3").WrapText = $true

# 4. J3 (Achieving % placeholder) now carries the percentage format that used to
#    belong to the old K3 cell.
$ws.Range("J3").NumberFormat = "0.00%"

# 5. Row 3's explicit 30pt height is no longer needed - let it size back to default.
$ws.Rows(3).AutoFit()

# 6. The merged header banner A1:K1 shrank to A1:J1 with the column removal; give the
#    new rightmost cell (J1) the closing right-hand border of that banner box.
$ws.Range("J1").Borders.Item(10).LineStyle = 1

# 7. Selection/merge bookkeeping to match the new A1:J1 banner.
$ws.Range("A1:J1").Select()

# 8. The hidden _FilterDatabase defined name still points at the old $M$2 edge;
#    bring it in to the new last column $L$2.
$wb.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$2:`$L`$2"
